# TC_148: rename the battery-load column headers to match the new
# implementation and grow columns S:T so the longer labels fit, then move
# the active selection onto the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Row 8 holds the actual column-header labels used by the new loading
# detail implementation.
$ws.Range("S8").Value = "Alarm Current(A)"
$ws.Range("T8").Value = "Standby Current(A)"

# Widen S:T to fit the new text, then select the data cell.
$ws.Columns.Item(19).ColumnWidth = 16.833333333333332
$ws.Columns.Item(20).ColumnWidth = 18.833333333333332
$ws.Range("S8:T8").Select() | Out-Null
